# Auto-assembled Excel COM-interop script
# Adds "Gradient Boosting" (GB) model columns to the describe, describe_relative
# and profit_stats sheets (one new column after each existing gNB/RF/DT/KNN model
# group on every sheet), and refreshes every score in all three tables to match
# the re-run cross-validation results that include the new GB model.

$wb = $excel.ActiveWorkbook

# ---- describe ----
$ws1 = $wb.Worksheets.Item("describe")
$ws1.Range("F1").EntireColumn.Insert()
$ws1.Range("K1").EntireColumn.Insert()

$ws1.Range("C1").Value = "FTR_RF"
$ws1.Range("D1").Value = "FTR_DT"
$ws1.Range("E1").Value = "FTR_KNN"
$ws1.Range("F1").Value = "FTR_GB"
$ws1.Range("G1").Value = "BTTS_gNB"
$ws1.Range("H1").Value = "BTTS_RF"
$ws1.Range("I1").Value = "BTTS_DT"
$ws1.Range("J1").Value = "BTTS_KNN"
$ws1.Range("K1").Value = "BTTS_GB"
$ws1.Range("L1").Value = "O/U2.5_gNB"
$ws1.Range("M1").Value = "O/U2.5_RF"
$ws1.Range("N1").Value = "O/U2.5_DT"
$ws1.Range("O1").Value = "O/U2.5_KNN"
$ws1.Range("P1").Value = "O/U2.5_GB"
$ws1.Range("B1").Copy() | Out-Null
$ws1.Range("C1:P1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws1.Range("B2").Value = 0.4974050632911393
$ws1.Range("C2").Value = 0.5082278481012659
$ws1.Range("D2").Value = 0.4633544303797468
$ws1.Range("E2").Value = 0.4681012658227848
$ws1.Range("F2").Value = 0.4963291139240507
$ws1.Range("G2").Value = 0.5965822784810126
$ws1.Range("H2").Value = 0.600379746835443
$ws1.Range("I2").Value = 0.6011392405063291
$ws1.Range("J2").Value = 0.5767088607594937
$ws1.Range("K2").Value = 0.6322784810126582
$ws1.Range("L2").Value = 0.6251898734177215
$ws1.Range("M2").Value = 0.6374683544303796
$ws1.Range("N2").Value = 0.6336075949367088
$ws1.Range("O2").Value = 0.6203164556962025
$ws1.Range("P2").Value = 0.6518987341772151

$ws1.Range("B3").Value = 0.05006998151802123
$ws1.Range("C3").Value = 0.05485762118443563
$ws1.Range("D3").Value = 0.05758093749176584
$ws1.Range("E3").Value = 0.05181123657631183
$ws1.Range("F3").Value = 0.05266866728244041
$ws1.Range("G3").Value = 0.04583020293839912
$ws1.Range("H3").Value = 0.04422439657662171
$ws1.Range("I3").Value = 0.05800020770818953
$ws1.Range("J3").Value = 0.04587111986991489
$ws1.Range("K3").Value = 0.04530789429936875
$ws1.Range("L3").Value = 0.04631047581207948
$ws1.Range("M3").Value = 0.04295970162544888
$ws1.Range("N3").Value = 0.04999595380899909
$ws1.Range("O3").Value = 0.04609625614387271
$ws1.Range("P3").Value = 0.04499134818491073

$ws1.Range("B4").Value = 0.3670886075949367
$ws1.Range("C4").Value = 0.3670886075949367
$ws1.Range("D4").Value = 0.2911392405063291
$ws1.Range("E4").Value = 0.3291139240506329
$ws1.Range("F4").Value = 0.3291139240506329
$ws1.Range("G4").Value = 0.4683544303797468
$ws1.Range("H4").Value = 0.4810126582278481
$ws1.Range("I4").Value = 0.4177215189873418
$ws1.Range("J4").Value = 0.4556962025316456
$ws1.Range("K4").Value = 0.5063291139240507
$ws1.Range("L4").Value = 0.5063291139240507
$ws1.Range("M4").Value = 0.5189873417721519
$ws1.Range("N4").Value = 0.4556962025316456
$ws1.Range("O4").Value = 0.4683544303797468
$ws1.Range("P4").Value = 0.5189873417721519

$ws1.Range("B5").Value = 0.4683544303797468
$ws1.Range("C5").Value = 0.4683544303797468
$ws1.Range("D5").Value = 0.4177215189873418
$ws1.Range("E5").Value = 0.4303797468354431
$ws1.Range("F5").Value = 0.4556962025316456
$ws1.Range("G5").Value = 0.569620253164557
$ws1.Range("H5").Value = 0.569620253164557
$ws1.Range("I5").Value = 0.569620253164557
$ws1.Range("J5").Value = 0.5443037974683544
$ws1.Range("K5").Value = 0.5949367088607594
$ws1.Range("L5").Value = 0.5949367088607594
$ws1.Range("M5").Value = 0.6075949367088608
$ws1.Range("N5").Value = 0.6075949367088608
$ws1.Range("O5").Value = 0.5949367088607594
$ws1.Range("P5").Value = 0.620253164556962

$ws1.Range("B6").Value = 0.4936708860759494
$ws1.Range("C6").Value = 0.5063291139240507
$ws1.Range("D6").Value = 0.4683544303797468
$ws1.Range("E6").Value = 0.4683544303797468
$ws1.Range("F6").Value = 0.4936708860759494
$ws1.Range("G6").Value = 0.6075949367088608
$ws1.Range("H6").Value = 0.5949367088607594
$ws1.Range("I6").Value = 0.6075949367088608
$ws1.Range("J6").Value = 0.5822784810126582
$ws1.Range("K6").Value = 0.6329113924050633
$ws1.Range("L6").Value = 0.620253164556962
$ws1.Range("M6").Value = 0.6329113924050633
$ws1.Range("N6").Value = 0.6329113924050633
$ws1.Range("O6").Value = 0.620253164556962
$ws1.Range("P6").Value = 0.6455696202531646

$ws1.Range("B7").Value = 0.5316455696202531
$ws1.Range("C7").Value = 0.5443037974683544
$ws1.Range("D7").Value = 0.5063291139240507
$ws1.Range("E7").Value = 0.5063291139240507
$ws1.Range("F7").Value = 0.5316455696202531
$ws1.Range("G7").Value = 0.6329113924050633
$ws1.Range("H7").Value = 0.6329113924050633
$ws1.Range("I7").Value = 0.6455696202531646
$ws1.Range("J7").Value = 0.6075949367088608
$ws1.Range("K7").Value = 0.6708860759493671
$ws1.Range("L7").Value = 0.6582278481012658
$ws1.Range("M7").Value = 0.6708860759493671
$ws1.Range("N7").Value = 0.6708860759493671
$ws1.Range("O7").Value = 0.6455696202531646
$ws1.Range("P7").Value = 0.6835443037974683

$ws1.Range("B8").Value = 0.620253164556962
$ws1.Range("C8").Value = 0.6708860759493671
$ws1.Range("D8").Value = 0.6329113924050633
$ws1.Range("E8").Value = 0.6329113924050633
$ws1.Range("F8").Value = 0.620253164556962
$ws1.Range("G8").Value = 0.6962025316455697
$ws1.Range("H8").Value = 0.7088607594936709
$ws1.Range("I8").Value = 0.7088607594936709
$ws1.Range("J8").Value = 0.7215189873417721
$ws1.Range("K8").Value = 0.759493670886076
$ws1.Range("L8").Value = 0.7468354430379747
$ws1.Range("M8").Value = 0.759493670886076
$ws1.Range("N8").Value = 0.7721518987341772
$ws1.Range("O8").Value = 0.7468354430379747
$ws1.Range("P8").Value = 0.759493670886076

# ---- describe_relative ----
$ws2 = $wb.Worksheets.Item("describe_relative")
$ws2.Range("F1").EntireColumn.Insert()
$ws2.Range("K1").EntireColumn.Insert()

$ws2.Range("C1").Value = "FTR_RF"
$ws2.Range("D1").Value = "FTR_DT"
$ws2.Range("E1").Value = "FTR_KNN"
$ws2.Range("F1").Value = "FTR_GB"
$ws2.Range("G1").Value = "BTTS_gNB"
$ws2.Range("H1").Value = "BTTS_RF"
$ws2.Range("I1").Value = "BTTS_DT"
$ws2.Range("J1").Value = "BTTS_KNN"
$ws2.Range("K1").Value = "BTTS_GB"
$ws2.Range("L1").Value = "O/U2.5_gNB"
$ws2.Range("M1").Value = "O/U2.5_RF"
$ws2.Range("N1").Value = "O/U2.5_DT"
$ws2.Range("O1").Value = "O/U2.5_KNN"
$ws2.Range("P1").Value = "O/U2.5_GB"
$ws2.Range("B1").Copy() | Out-Null
$ws2.Range("C1:P1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws2.Range("B2").Value = 49.2
$ws2.Range("C2").Value = 52.5
$ws2.Range("D2").Value = 39
$ws2.Range("E2").Value = 40.4
$ws2.Range("F2").Value = 48.9
$ws2.Range("G2").Value = 19.3
$ws2.Range("H2").Value = 20.1
$ws2.Range("I2").Value = 20.2
$ws2.Range("J2").Value = 15.3
$ws2.Range("K2").Value = 26.5
$ws2.Range("L2").Value = 25
$ws2.Range("M2").Value = 27.5
$ws2.Range("N2").Value = 26.7
$ws2.Range("O2").Value = 24.1
$ws2.Range("P2").Value = 30.4

$ws2.Range("B3").Value = -85
$ws2.Range("C3").Value = -83.5
$ws2.Range("D3").Value = -82.7
$ws2.Range("E3").Value = -84.5
$ws2.Range("F3").Value = -84.2
$ws2.Range("G3").Value = -90.8
$ws2.Range("H3").Value = -91.2
$ws2.Range("I3").Value = -88.40000000000001
$ws2.Range("J3").Value = -90.8
$ws2.Range("K3").Value = -90.90000000000001
$ws2.Range("L3").Value = -90.7
$ws2.Range("M3").Value = -91.40000000000001
$ws2.Range("N3").Value = -90
$ws2.Range("O3").Value = -90.8
$ws2.Range("P3").Value = -91

$ws2.Range("B4").Value = 10.1
$ws2.Range("C4").Value = 10.1
$ws2.Range("D4").Value = -12.7
$ws2.Range("E4").Value = -1.3
$ws2.Range("F4").Value = -1.3
$ws2.Range("G4").Value = -6.3
$ws2.Range("H4").Value = -3.8
$ws2.Range("I4").Value = -16.5
$ws2.Range("J4").Value = -8.9
$ws2.Range("K4").Value = 1.3
$ws2.Range("L4").Value = 1.3
$ws2.Range("M4").Value = 3.8
$ws2.Range("N4").Value = -8.9
$ws2.Range("O4").Value = -6.3
$ws2.Range("P4").Value = 3.8

$ws2.Range("B5").Value = 40.5
$ws2.Range("C5").Value = 40.5
$ws2.Range("D5").Value = 25.3
$ws2.Range("E5").Value = 29.1
$ws2.Range("F5").Value = 36.7
$ws2.Range("G5").Value = 13.9
$ws2.Range("H5").Value = 13.9
$ws2.Range("I5").Value = 13.9
$ws2.Range("J5").Value = 8.9
$ws2.Range("K5").Value = 19
$ws2.Range("L5").Value = 19
$ws2.Range("M5").Value = 21.5
$ws2.Range("N5").Value = 21.5
$ws2.Range("O5").Value = 19
$ws2.Range("P5").Value = 24.1

$ws2.Range("B6").Value = 48.1
$ws2.Range("C6").Value = 51.9
$ws2.Range("D6").Value = 40.5
$ws2.Range("E6").Value = 40.5
$ws2.Range("F6").Value = 48.1
$ws2.Range("G6").Value = 21.5
$ws2.Range("H6").Value = 19
$ws2.Range("I6").Value = 21.5
$ws2.Range("J6").Value = 16.5
$ws2.Range("K6").Value = 26.6
$ws2.Range("L6").Value = 24.1
$ws2.Range("M6").Value = 26.6
$ws2.Range("N6").Value = 26.6
$ws2.Range("O6").Value = 24.1
$ws2.Range("P6").Value = 29.1

$ws2.Range("B7").Value = 59.5
$ws2.Range("C7").Value = 63.3
$ws2.Range("D7").Value = 51.9
$ws2.Range("E7").Value = 51.9
$ws2.Range("F7").Value = 59.5
$ws2.Range("G7").Value = 26.6
$ws2.Range("H7").Value = 26.6
$ws2.Range("I7").Value = 29.1
$ws2.Range("J7").Value = 21.5
$ws2.Range("K7").Value = 34.2
$ws2.Range("L7").Value = 31.6
$ws2.Range("M7").Value = 34.2
$ws2.Range("N7").Value = 34.2
$ws2.Range("O7").Value = 29.1
$ws2.Range("P7").Value = 36.7

$ws2.Range("B8").Value = 86.09999999999999
$ws2.Range("C8").Value = 101.3
$ws2.Range("D8").Value = 89.90000000000001
$ws2.Range("E8").Value = 89.90000000000001
$ws2.Range("F8").Value = 86.09999999999999
$ws2.Range("G8").Value = 39.2
$ws2.Range("H8").Value = 41.8
$ws2.Range("I8").Value = 41.8
$ws2.Range("J8").Value = 44.3
$ws2.Range("K8").Value = 51.9
$ws2.Range("L8").Value = 49.4
$ws2.Range("M8").Value = 51.9
$ws2.Range("N8").Value = 54.4
$ws2.Range("O8").Value = 49.4
$ws2.Range("P8").Value = 51.9

# ---- profit_stats ----
$ws3 = $wb.Worksheets.Item("profit_stats")
$ws3.Range("F1").EntireColumn.Insert()

$ws3.Range("C1").Value = "FTR_RF"
$ws3.Range("D1").Value = "FTR_DT"
$ws3.Range("E1").Value = "FTR_KNN"
$ws3.Range("F1").Value = "FTR_GB"
$ws3.Range("G1").Value = "O/U2.5_gNB"
$ws3.Range("H1").Value = "O/U2.5_RF"
$ws3.Range("I1").Value = "O/U2.5_DT"
$ws3.Range("J1").Value = "O/U2.5_KNN"
$ws3.Range("K1").Value = "O/U2.5_GB"
$ws3.Range("B1").Copy() | Out-Null
$ws3.Range("C1:K1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws3.Range("B2").Value = 300
$ws3.Range("C2").Value = 300
$ws3.Range("D2").Value = 300
$ws3.Range("E2").Value = 300
$ws3.Range("F2").Value = 300
$ws3.Range("G2").Value = 300
$ws3.Range("H2").Value = 300
$ws3.Range("I2").Value = 300
$ws3.Range("J2").Value = 300
$ws3.Range("K2").Value = 300

$ws3.Range("B3").Value = -8.33549258465899
$ws3.Range("C3").Value = -19.67085744922444
$ws3.Range("D3").Value = -17.19022348328681
$ws3.Range("E3").Value = -21.98583159884765
$ws3.Range("F3").Value = -25.47855093993547
$ws3.Range("G3").Value = 4.475051963100621
$ws3.Range("H3").Value = 1.183489052441255
$ws3.Range("I3").Value = 4.440264506146576
$ws3.Range("J3").Value = 1.037060247171797
$ws3.Range("K3").Value = 5.182611128559667

$ws3.Range("B4").Value = 10.65040341048969
$ws3.Range("C4").Value = 13.49452794585846
$ws3.Range("D4").Value = 13.05067505494932
$ws3.Range("E4").Value = 12.13864442288408
$ws3.Range("F4").Value = 14.65410307287836
$ws3.Range("G4").Value = 10.34964384037962
$ws3.Range("H4").Value = 9.684320781986543
$ws3.Range("I4").Value = 10.41363468730878
$ws3.Range("J4").Value = 8.958648321727878
$ws3.Range("K4").Value = 10.31973280993218

$ws3.Range("B5").Value = -38.22271468001645
$ws3.Range("C5").Value = -62.576647032816
$ws3.Range("D5").Value = -60.94657844214715
$ws3.Range("E5").Value = -56.37842191290093
$ws3.Range("F5").Value = -72.41503743072639
$ws3.Range("G5").Value = -27.66962832675908
$ws3.Range("H5").Value = -28.23598459712059
$ws3.Range("I5").Value = -29.20245926209875
$ws3.Range("J5").Value = -26.45552385504812
$ws3.Range("K5").Value = -27.66962832675908

$ws3.Range("B6").Value = -15.17254330629996
$ws3.Range("C6").Value = -28.02162266470328
$ws3.Range("D6").Value = -26.4406484976103
$ws3.Range("E6").Value = -31.14011213103829
$ws3.Range("F6").Value = -35.41716570695253
$ws3.Range("G6").Value = -2.161147282523468
$ws3.Range("H6").Value = -5.471963685094597
$ws3.Range("I6").Value = -2.123155510459363
$ws3.Range("J6").Value = -4.504321018487849
$ws3.Range("K6").Value = -1.410253170161539

$ws3.Range("B7").Value = -7.955298980566197
$ws3.Range("C7").Value = -18.94198887514706
$ws3.Range("D7").Value = -16.48022981490216
$ws3.Range("E7").Value = -21.27224949026942
$ws3.Range("F7").Value = -24.67420369192573
$ws3.Range("G7").Value = 4.558164148300445
$ws3.Range("H7").Value = 1.268508260673014
$ws3.Range("I7").Value = 4.126965459878039
$ws3.Range("J7").Value = 1.359808983881096
$ws3.Range("K7").Value = 5.175290314862615

$ws3.Range("B8").Value = -0.4244789575371394
$ws3.Range("C8").Value = -10.50899140920108
$ws3.Range("D8").Value = -7.838569618802591
$ws3.Range("E8").Value = -13.10107099211462
$ws3.Range("F8").Value = -15.18324506394002
$ws3.Range("G8").Value = 11.54392964133005
$ws3.Range("H8").Value = 7.622170647822144
$ws3.Range("I8").Value = 11.60531976375231
$ws3.Range("J8").Value = 7.186151140689343
$ws3.Range("K8").Value = 11.85703828579945

$ws3.Range("B9").Value = 18.24998841398354
$ws3.Range("C9").Value = 15.55098470067806
$ws3.Range("D9").Value = 15.1445810611478
$ws3.Range("E9").Value = 8.717009939740446
$ws3.Range("F9").Value = 11.45456757226248
$ws3.Range("G9").Value = 32.35565549616786
$ws3.Range("H9").Value = 31.7102343601391
$ws3.Range("I9").Value = 32.39580169034406
$ws3.Range("J9").Value = 24.02101019025249
$ws3.Range("K9").Value = 41.20693309117006

